$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 10; existing rows 10-21 shift down to 11-22
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with data
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = "Vega Monumental Concepción"
$ws.Range("C10").Value = "Bíobío"
$ws.Range("D10").Value = 45280
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100103
$ws.Range("H10").Value = "Frutos de hueso (carozo)"
$ws.Range("I10").Value = 100103003
$ws.Range("J10").Value = "Damasco"
$ws.Range("K10").Value = "Castle Brite"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 180
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 13000
$ws.Range("P10").Value = 12556
$ws.Range("Q10").Value = "`$/bandeja 10 kilos"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 1256
$ws.Range("T10").Value = 10
